$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 10.51982755881343
$ws.Cells.Item(2, 3).Value = 6.428754427464815
$ws.Cells.Item(2, 4).Value = 6.329196796453393
$ws.Cells.Item(2, 6).Value = 31.64595836692241
$ws.Cells.Item(2, 7).Value = 3.676362259101269
$ws.Cells.Item(2, 9).Value = 25.26740506135005
$ws.Cells.Item(2, 11).Value = 11.04650988020718
# Row 3
$ws.Cells.Item(3, 2).Value = 10.26763139553136
$ws.Cells.Item(3, 3).Value = 6.231916815297287
$ws.Cells.Item(3, 4).Value = 6.324345970649012
$ws.Cells.Item(3, 6).Value = 31.64030690474888
$ws.Cells.Item(3, 7).Value = 3.678906672085041
$ws.Cells.Item(3, 9).Value = 25.34066460789988
$ws.Cells.Item(3, 11).Value = 10.88088821425215
# Row 4
$ws.Cells.Item(4, 2).Value = 10.11218708096266
$ws.Cells.Item(4, 3).Value = 6.109775615690326
$ws.Cells.Item(4, 4).Value = 6.321357007615719
$ws.Cells.Item(4, 6).Value = 31.64561919729205
$ws.Cells.Item(4, 7).Value = 3.680550290595583
$ws.Cells.Item(4, 9).Value = 25.39101147492261
$ws.Cells.Item(4, 11).Value = 10.7806025609717
# Row 5
$ws.Cells.Item(5, 2).Value = 10.04879400688309
$ws.Cells.Item(5, 3).Value = 6.059766855820411
$ws.Cells.Item(5, 4).Value = 6.320136782924838
$ws.Cells.Item(5, 6).Value = 31.64998832856073
$ws.Cells.Item(5, 7).Value = 3.68124060062393
$ws.Cells.Item(5, 9).Value = 25.41287256691028
$ws.Cells.Item(5, 11).Value = 10.74014237511744
# Row 6
$ws.Cells.Item(6, 2).Value = 10.038267750677
$ws.Cells.Item(6, 3).Value = 6.051451383162525
$ws.Cells.Item(6, 4).Value = 6.319934046163603
$ws.Cells.Item(6, 6).Value = 31.65084680546236
$ws.Cells.Item(6, 7).Value = 3.681356467589344
$ws.Cells.Item(6, 9).Value = 25.41658363632526
$ws.Cells.Item(6, 11).Value = 10.73345012409743
# Row 7
$ws.Cells.Item(7, 2).Value = 10.11133219869517
$ws.Cells.Item(7, 3).Value = 6.10910201340624
$ws.Cells.Item(7, 4).Value = 6.321340559551531
$ws.Cells.Item(7, 6).Value = 31.64566920230948
$ws.Cells.Item(7, 7).Value = 3.680559517173541
$ws.Cells.Item(7, 9).Value = 25.391300864585
$ws.Cells.Item(7, 11).Value = 10.78005518263316
# Row 8
$ws.Cells.Item(8, 2).Value = 10.43305342087965
$ws.Cells.Item(8, 3).Value = 6.361201532045689
$ws.Cells.Item(8, 4).Value = 6.327526685645475
$ws.Cells.Item(8, 6).Value = 31.64218501050575
$ws.Cells.Item(8, 7).Value = 3.677222732405064
$ws.Cells.Item(8, 9).Value = 25.2915485642852
$ws.Cells.Item(8, 11).Value = 10.98914365453551
# Row 9
$ws.Cells.Item(9, 2).Value = 11.05497882112514
$ws.Cells.Item(9, 3).Value = 6.84185048399214
$ws.Cells.Item(9, 4).Value = 6.339555049841927
$ws.Cells.Item(9, 6).Value = 31.70514760522741
$ws.Cells.Item(9, 7).Value = 3.671321546243411
$ws.Cells.Item(9, 9).Value = 25.13871413953472
$ws.Cells.Item(9, 11).Value = 11.40797543419919
# Row 10
$ws.Cells.Item(10, 2).Value = 11.50108465648057
$ws.Cells.Item(10, 3).Value = 7.182299191858804
$ws.Cells.Item(10, 4).Value = 6.348305396056821
$ws.Cells.Item(10, 6).Value = 31.79400513297488
$ws.Cells.Item(10, 7).Value = 3.667373059166811
$ws.Cells.Item(10, 9).Value = 25.0527842603239
$ws.Cells.Item(10, 11).Value = 11.71793149723866
# Row 11
$ws.Cells.Item(11, 2).Value = 11.70065431772022
$ws.Cells.Item(11, 3).Value = 7.333634527782398
$ws.Cells.Item(11, 4).Value = 6.352262317100391
$ws.Cells.Item(11, 6).Value = 31.8436534969533
$ws.Cells.Item(11, 7).Value = 3.66565990821524
$ws.Cells.Item(11, 9).Value = 25.01947213089844
$ws.Cells.Item(11, 11).Value = 11.85877739345491
# Row 12
$ws.Cells.Item(12, 2).Value = 11.77566459461868
$ws.Cells.Item(12, 3).Value = 7.390375315878686
$ws.Cells.Item(12, 4).Value = 6.353756880301079
$ws.Cells.Item(12, 6).Value = 31.86377565781909
$ws.Cells.Item(12, 7).Value = 3.665023051481125
$ws.Cells.Item(12, 9).Value = 25.00769287154812
$ws.Cells.Item(12, 11).Value = 11.91203761943228
# Row 13
$ws.Cells.Item(13, 2).Value = 11.75953598732747
$ws.Cells.Item(13, 3).Value = 7.378181234227342
$ws.Cells.Item(13, 4).Value = 6.353435176858965
$ws.Cells.Item(13, 6).Value = 31.85938332526253
$ws.Cells.Item(13, 7).Value = 3.66515968278578
$ws.Cells.Item(13, 9).Value = 25.0101925254925
$ws.Cells.Item(13, 11).Value = 11.90057123672583
# Row 14
$ws.Cells.Item(14, 2).Value = 11.70683715770235
$ws.Cells.Item(14, 3).Value = 7.338314294027727
$ws.Cells.Item(14, 4).Value = 6.352385355075654
$ws.Cells.Item(14, 6).Value = 31.84528251349627
$ws.Cells.Item(14, 7).Value = 3.665607275974858
$ws.Cells.Item(14, 9).Value = 25.01848627812761
$ws.Cells.Item(14, 11).Value = 11.86316094468602
# Row 15
$ws.Cells.Item(15, 2).Value = 11.67448213534095
$ws.Cells.Item(15, 3).Value = 7.313819229289674
$ws.Cells.Item(15, 4).Value = 6.351741796463695
$ws.Cells.Item(15, 6).Value = 31.83681725477076
$ws.Cells.Item(15, 7).Value = 3.665882984483818
$ws.Cells.Item(15, 9).Value = 25.0236753540538
$ws.Cells.Item(15, 11).Value = 11.84023474893656
# Row 16
$ws.Cells.Item(16, 2).Value = 11.48796773762543
$ws.Cells.Item(16, 3).Value = 7.172332807594859
$ws.Cells.Item(16, 4).Value = 6.348046288547449
$ws.Cells.Item(16, 6).Value = 31.79094574036911
$ws.Cells.Item(16, 7).Value = 3.667486683012253
$ws.Cells.Item(16, 9).Value = 25.05507793805079
$ws.Cells.Item(16, 11).Value = 11.70871921922979
# Row 17
$ws.Cells.Item(17, 2).Value = 11.37262737932519
$ws.Cells.Item(17, 3).Value = 7.084586464384433
$ws.Cells.Item(17, 4).Value = 6.345772843026849
$ws.Cells.Item(17, 6).Value = 31.76516495739729
$ws.Cells.Item(17, 7).Value = 3.668491721809694
$ws.Cells.Item(17, 9).Value = 25.07582549134295
$ws.Cells.Item(17, 11).Value = 11.62796095892332
# Row 18
$ws.Cells.Item(18, 2).Value = 11.30597412735789
$ws.Cells.Item(18, 3).Value = 7.033787301391607
$ws.Cells.Item(18, 4).Value = 6.344463016173036
$ws.Cells.Item(18, 6).Value = 31.75120542244497
$ws.Cells.Item(18, 7).Value = 3.669077613051409
$ws.Cells.Item(18, 9).Value = 25.08830251377017
$ws.Cells.Item(18, 11).Value = 11.58149944256222
# Row 19
$ws.Cells.Item(19, 2).Value = 11.28335529104953
$ws.Cells.Item(19, 3).Value = 7.016532804100796
$ws.Cells.Item(19, 4).Value = 6.344019166681231
$ws.Cells.Item(19, 6).Value = 31.7466283309125
$ws.Cells.Item(19, 7).Value = 3.669277330726345
$ws.Cells.Item(19, 9).Value = 25.09262023959744
$ws.Cells.Item(19, 11).Value = 11.56576803401389
# Row 20
$ws.Cells.Item(20, 2).Value = 11.38493848972054
$ws.Cells.Item(20, 3).Value = 7.0939617837738
$ws.Cells.Item(20, 4).Value = 6.346015086156282
$ws.Cells.Item(20, 6).Value = 31.76781947214052
$ws.Cells.Item(20, 7).Value = 3.668383924889442
$ws.Cells.Item(20, 9).Value = 25.07356058739417
$ws.Cells.Item(20, 11).Value = 11.63655939211146
# Row 21
$ws.Cells.Item(21, 2).Value = 11.7223319402162
$ws.Cells.Item(21, 3).Value = 7.350039992063458
$ws.Cells.Item(21, 4).Value = 6.352693820740376
$ws.Cells.Item(21, 6).Value = 31.84938845291328
$ws.Cells.Item(21, 7).Value = 3.665475485193052
$ws.Cells.Item(21, 9).Value = 25.01602749686042
$ws.Cells.Item(21, 11).Value = 11.87415170332153
# Row 22
$ws.Cells.Item(22, 2).Value = 11.9395277164818
$ws.Cells.Item(22, 3).Value = 7.514076402572299
$ws.Cells.Item(22, 4).Value = 6.3570361993983
$ws.Cells.Item(22, 6).Value = 31.91039689422186
$ws.Cells.Item(22, 7).Value = 3.66364384674483
$ws.Cells.Item(22, 9).Value = 24.98329711998346
$ws.Cells.Item(22, 11).Value = 12.02897134267525
# Row 23
$ws.Cells.Item(23, 2).Value = 11.82393348146504
$ws.Cells.Item(23, 3).Value = 7.426849033381141
$ws.Cells.Item(23, 4).Value = 6.35472079611518
$ws.Cells.Item(23, 6).Value = 31.87713337361264
$ws.Cells.Item(23, 7).Value = 3.664615116367981
$ws.Cells.Item(23, 9).Value = 25.00031883032476
$ws.Cells.Item(23, 11).Value = 11.94640030997706
# Row 24
$ws.Cells.Item(24, 2).Value = 11.37937369762242
$ws.Cells.Item(24, 3).Value = 7.089724295750058
$ws.Cells.Item(24, 4).Value = 6.345905576674474
$ws.Cells.Item(24, 6).Value = 31.76661667990992
$ws.Cells.Item(24, 7).Value = 3.668432634699418
$ws.Cells.Item(24, 9).Value = 25.07458284069258
$ws.Cells.Item(24, 11).Value = 11.63267213768816
# Row 25
$ws.Cells.Item(25, 2).Value = 10.88826498802085
$ws.Cells.Item(25, 3).Value = 6.713761026642088
$ws.Cells.Item(25, 4).Value = 6.336314590502766
$ws.Cells.Item(25, 6).Value = 31.68063048897474
$ws.Cells.Item(25, 7).Value = 3.672849674981649
$ws.Cells.Item(25, 9).Value = 25.26740506135005
$ws.Cells.Item(25, 11).Value = 11.29405300026637
